$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2021" column (J) to the right of the existing "2020" column (I),
# reusing the same cell formatting (number formats, borders, etc.) that
# column I already has for the header row and each data row.
$ws.Range("I4:I14").Copy() | Out-Null
$ws.Range("J4:J14").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Header
$ws.Range("J4").Value = 2021

# Data values for 2021
$ws.Range("J5").Value = 24.4
$ws.Range("J6").Value = 45.7
$ws.Range("J7").Value = 38
$ws.Range("J8").Value = 51.3
$ws.Range("J9").Value = 51.5
$ws.Range("J10").Value = 13
$ws.Range("J11").Value = 36.4
$ws.Range("J12").Value = 27
$ws.Range("J13").Value = 2.7
$ws.Range("J14").Value = 40.4

# Tweak row 3's height slightly
$ws.Rows.Item(3).RowHeight = 13.5

# Move/restore the active selection
$ws.Range("K18").Select() | Out-Null
